$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $text) {
    $cell = $table.Cell($row, 1)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $text
}

# Single-value cells whose numbers changed
Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "40"
Set-CellText $t 6 "0.00056"
Set-CellText $t 7 "0.00025"
Set-CellText $t 8 "0.00005"
Set-CellText $t 9 "0.00047"
Set-CellText $t 10 "0.00052"
Set-CellText $t 11 "0.00055"
Set-CellText $t 12 "0.00996"

# Cells that previously held a tab-separated run of many values are
# collapsed down to a single value
Set-CellText $t 44 "99.95"
Set-CellText $t 45 "0.01"
Set-CellText $t 46 "19"
